# 141: 31/12 10:23 LP1912+6203+6173
# Updates the "Última actualización" timestamp and "Total filas" counts,
# and appends newly-scraped rows to the LP1912 and 6203-6173 sheets.

$wb = $excel.ActiveWorkbook

$newTimestamp = "Última actualización: 31/12/2025 07:23:07"

# ---------------------------------------------------------------------
# Sheet "LP1912": refresh header info + append rows 704-716
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = $newTimestamp
$ws1.Range("A3").Value = "Total filas: 715"

$sheet1Rows = @(
    @("07:22:56", "07:28", "14_ABASTO",     6,  "LP1912", "31/12/2025"),
    @("07:22:56", "07:33", "23_HERNANDEZ",  11, "LP1912", "31/12/2025"),
    @("07:22:56", "07:35", "17X38_ROMERO",  13, "LP1912", "31/12/2025"),
    @("07:22:56", "07:36", "27_EL RETIRO",  14, "LP1912", "31/12/2025"),
    @("07:22:56", "07:51", "15_ABASTO",     29, "LP1912", "31/12/2025"),
    @("07:22:56", "07:58", "23_HERNANDEZ",  36, "LP1912", "31/12/2025"),
    @("07:22:56", "08:01", "16_SANTA ANA",  39, "LP1912", "31/12/2025"),
    @("07:22:56", "08:03", "17_ROMERO",     41, "LP1912", "31/12/2025"),
    @("07:22:56", "08:13", "10_OLMOS",      51, "LP1912", "31/12/2025"),
    @("07:22:56", "08:15", "17_ROMERO",     53, "LP1912", "31/12/2025"),
    @("07:22:56", "08:29", "14_ABASTO",     67, "LP1912", "31/12/2025"),
    @("07:22:56", "08:43", "10_OLMOS",      81, "LP1912", "31/12/2025"),
    @("07:22:56", "08:49", "16_SANTA ANA",  87, "LP1912", "31/12/2025")
)

$row = 704
foreach ($r in $sheet1Rows) {
    $ws1.Cells.Item($row, 2).Value = $r[0]
    $ws1.Cells.Item($row, 3).Value = $r[1]
    $ws1.Cells.Item($row, 4).Value = $r[2]
    $ws1.Cells.Item($row, 5).Value = $r[3]
    $ws1.Cells.Item($row, 6).Value = $r[4]
    $ws1.Cells.Item($row, 7).Value = $r[5]
    $row = $row + 1
}

# ---------------------------------------------------------------------
# Sheet "LP1912-215": refresh header info only (no new rows)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $newTimestamp

# ---------------------------------------------------------------------
# Sheet "6203-6173": refresh header info + append rows 87-89
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $newTimestamp
$ws3.Range("A3").Value = "Total filas: 88"

$sheet3Rows = @(
    @("31/12/2025", "07:23:06", "07:30", "215A_LA PLATA", 7,  "L6173"),
    @("31/12/2025", "07:23:06", "08:10", "215A_LA PLATA", 47, "L6173"),
    @("31/12/2025", "07:23:01", "08:36", "215C_LA PLATA", 73, "L6203")
)

$row = 87
foreach ($r in $sheet3Rows) {
    $ws3.Cells.Item($row, 2).Value = $r[0]
    $ws3.Cells.Item($row, 3).Value = $r[1]
    $ws3.Cells.Item($row, 4).Value = $r[2]
    $ws3.Cells.Item($row, 5).Value = $r[3]
    $ws3.Cells.Item($row, 6).Value = $r[4]
    $ws3.Cells.Item($row, 7).Value = $r[5]
    $row = $row + 1
}

Write-Output "edit applied"
